# Automation HUB process: refresh the Yearly-Report data rows.
# Row 1-2 are left untouched; rows 3-11 are (re)written with the
# report's current record set (3 brand-new invoice rows mixed in,
# table grown from 8 to 11 data rows).
#
# All source data is plain text (invoice/amount columns included),
# so every cell is forced to Text format before the write and the
# formatting is cleared again afterwards - this keeps the values as
# shared-string text (matching the original table) instead of letting
# Excel auto-coerce numeric- or date-looking strings into real
# numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(3,  "365255", "Professional Services",  "218977", "43795.4", "262772",  "EUR", "2017-06-21"),
    @(4,  "611337", "IT Support",             "274278", "54855.6", "329134",  "USD", "2017-05-10"),
    @(5,  "303148", "Concierge Services",     "166873", "33374.6", "200248",  "EUR", "2017-05-07"),
    @(6,  "365255", "Professional Services",  "218977", "43795.4", "262772",  "EUR", "2017-06-21"),
    @(7,  "762404", "Professional Services",  "150414", "30082.8", "180497",  "RON", "2017-07-21"),
    @(8,  "259639", "Concierge Services",     "257124", "51424.8", "308549",  "RON", "2017-08-24"),
    @(9,  "550727", "Professional Services",  "158243", "31648.6", "189892",  "RON", "2017-10-07"),
    @(10, "863559", "Beverages and Catering",  "78842", "15768.4", "94610.4", "EUR", "2017-10-01"),
    @(11, "830889", "IT Support",             "248649", "49729.8", "298379",  "USD", "2017-11-27")
)

$fullRange = $ws.Range("A3:G11")
$fullRange.NumberFormat = "@"

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Range("G$r").Value = $row[7]
}

$fullRange.ClearFormats()
